$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Saturday hours for the week commencing 43171 (row 9)
$ws.Range("G9").Value = 6.75

# Update the selected cell/range as recorded in the saved view
$ws.Range("M12").Select()
